$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1954397394136808
$ws.Range("C2").Value = 0.5374592833876222
$ws.Range("J2").Value = 0.02605863192182411
$ws.Range("P2").Value = 0.1628664495114006
$ws.Range("S2").Value = 0.07817589576547231
$ws.Range("B3").Value = 0.01162790697674419
$ws.Range("C3").Value = 0.01162790697674419
$ws.Range("J3").Value = 0.03488372093023256
$ws.Range("P3").Value = 0.7267441860465116
$ws.Range("S3").Value = 0.2151162790697674
$ws.Range("J4").Value = 0.08333333333333333
$ws.Range("O4").Value = 0.02083333333333333
$ws.Range("P4").Value = 0.7291666666666666
$ws.Range("B6").Value = 0.05194805194805195
$ws.Range("F6").Value = 0.03463203463203463
$ws.Range("J6").Value = 0.3073593073593073
$ws.Range("O6").Value = 0.02164502164502164
$ws.Range("Q6").Value = 0.1731601731601732
$ws.Range("R6").Value = 0.06926406926406926
$ws.Range("S6").Value = 0.341991341991342
$ws.Range("B7").Value = 0.1063829787234043
$ws.Range("D7").Value = 0.04964539007092199
$ws.Range("F7").Value = 0.05673758865248227
$ws.Range("J7").Value = 0.1063829787234043
$ws.Range("O7").Value = 0.007092198581560284
$ws.Range("Q7").Value = 0.2269503546099291
$ws.Range("R7").Value = 0.07801418439716312
$ws.Range("S7").Value = 0.3687943262411347
$ws.Range("B8").Value = 0.08024691358024691
$ws.Range("D8").Value = 0.0308641975308642
$ws.Range("F8").Value = 0.07407407407407407
$ws.Range("J8").Value = 0.102880658436214
$ws.Range("O8").Value = 0.006172839506172839
$ws.Range("Q8").Value = 0.2098765432098765
$ws.Range("R8").Value = 0.06995884773662552
$ws.Range("S8").Value = 0.4259259259259259
$ws.Range("B9").Value = 0.1072961373390558
$ws.Range("D9").Value = 0.02145922746781116
$ws.Range("F9").Value = 0.05150214592274678
$ws.Range("J9").Value = 0.09012875536480687
$ws.Range("O9").Value = 0.0128755364806867
$ws.Range("Q9").Value = 0.2060085836909871
$ws.Range("R9").Value = 0.07296137339055794
$ws.Range("S9").Value = 0.4377682403433477
$ws.Range("B10").Value = 0.1226492232215863
$ws.Range("D10").Value = 0.01635322976287817
$ws.Range("F10").Value = 0.0678659035159444
$ws.Range("J10").Value = 0.1062959934587081
$ws.Range("O10").Value = 0.01798855273916599
$ws.Range("Q10").Value = 0.2150449713818479
$ws.Range("R10").Value = 0.07277187244480784
$ws.Range("S10").Value = 0.3810302534750613
$ws.Range("G11").Value = 0.1617647058823529
$ws.Range("J11").Value = 0.07843137254901961
$ws.Range("K11").Value = 0.196078431372549
$ws.Range("L11").Value = 0.5343137254901961
$ws.Range("S11").Value = 0.02941176470588235
$ws.Range("G12").Value = 0.7946428571428571
$ws.Range("J12").Value = 0.1607142857142857
$ws.Range("L12").Value = 0.03571428571428571
$ws.Range("S12").Value = 0.008928571428571428
$ws.Range("G13").Value = 0.6097560975609756
$ws.Range("J13").Value = 0.3414634146341464
$ws.Range("S13").Value = 0.04878048780487805
$ws.Range("F15").Value = 0.02145922746781116
$ws.Range("H15").Value = 0.1630901287553648
$ws.Range("I15").Value = 0.07296137339055794
$ws.Range("J15").Value = 0.3605150214592275
$ws.Range("K15").Value = 0.05579399141630902
$ws.Range("M15").Value = 0.01716738197424893
$ws.Range("O15").Value = 0.06437768240343347
$ws.Range("S15").Value = 0.2446351931330472
$ws.Range("F16").Value = 0.03431372549019608
$ws.Range("H16").Value = 0.1911764705882353
$ws.Range("I16").Value = 0.07352941176470588
$ws.Range("J16").Value = 0.4705882352941176
$ws.Range("K16").Value = 0.08333333333333333
$ws.Range("M16").Value = 0.009803921568627451
$ws.Range("O16").Value = 0.04411764705882353
$ws.Range("S16").Value = 0.09313725490196079
$ws.Range("F17").Value = 0.01646090534979424
$ws.Range("H17").Value = 0.2160493827160494
$ws.Range("I17").Value = 0.1172839506172839
$ws.Range("J17").Value = 0.3909465020576132
$ws.Range("K17").Value = 0.06790123456790123
$ws.Range("M17").Value = 0.01646090534979424
$ws.Range("N17").Value = 0.00205761316872428
$ws.Range("O17").Value = 0.05349794238683128
$ws.Range("S17").Value = 0.1193415637860082
$ws.Range("F18").Value = 0.01764705882352941
$ws.Range("H18").Value = 0.1647058823529412
$ws.Range("I18").Value = 0.09411764705882353
$ws.Range("J18").Value = 0.4176470588235294
$ws.Range("K18").Value = 0.04705882352941176
$ws.Range("O18").Value = 0.1
$ws.Range("S18").Value = 0.1588235294117647
$ws.Range("F19").Value = 0.02782071097372488
$ws.Range("H19").Value = 0.2187017001545595
$ws.Range("I19").Value = 0.1004636785162288
$ws.Range("J19").Value = 0.3454404945904173
$ws.Range("K19").Value = 0.0687789799072643
$ws.Range("M19").Value = 0.02163833075734158
$ws.Range("N19").Value = 0.0007727975270479134
$ws.Range("O19").Value = 0.08114374034003091
$ws.Range("S19").Value = 0.1352395672333848
